$wb = $excel.ActiveWorkbook
$ws15 = $wb.Worksheets.Item("15")
Write-Host $ws15.Range("A2").Style
Write-Host $ws15.Range("A3").Style
Write-Host $ws15.Range("B3").Style
Write-Host $ws15.Range("A1").Style
